# Update countries & provincias Spain
# Applies the data refresh that shifted the "Siria" row up (new stats),
# pushing "Angola" and "Polinesia Francesa" down by one row, and updates
# several other countries' daily COVID-19 figures plus the "last updated"
# timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp header in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 12:05"

# Belgica (row 21): Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$ws.Range("B21").Value = 56810
$ws.Range("C21").Value = 299
$ws.Range("D21").Value = 15155
$ws.Range("E21").Value = 32418
$ws.Range("G21").Value = 25
$ws.Range("H21").Value = 9237

# Rumania (row 40)
$ws.Range("E40").Value = 5765
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 1170

# Israel (row 41)
$ws.Range("B41").Value = 16712
$ws.Range("C41").Value = 22
$ws.Range("D41").Value = 14085
$ws.Range("E41").Value = 2348

# Albania (row 109)
$ws.Range("B109").Value = 989
$ws.Range("C109").Value = 8
$ws.Range("D109").Value = 783
$ws.Range("E109").Value = 175

# Siria now occupies row 178 with refreshed figures, Angola and Polinesia
# Francesa are pushed down to rows 179 and 180 respectively (their own
# figures are unchanged, only their row position shifts).
$ws.Range("A178").Value = "Siria"
$ws.Range("B178").Value = 70
$ws.Range("C178").Value = 11
$ws.Range("D178").Value = 37
$ws.Range("E178").Value = 29
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 4

$ws.Range("A179").Value = "Angola"
$ws.Range("B179").Value = 60
$ws.Range("C179").Value = 0
$ws.Range("D179").Value = 17
$ws.Range("E179").Value = 40
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 3

$ws.Range("A180").Value = "Polinesia Francesa"
$ws.Range("B180").Value = 60
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 60
$ws.Range("E180").Value = 0
$ws.Range("F180").Value = 0
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 0

$wb.Save()
